$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Update existing CasesTab row (row 2) query text: ["txt"] -> ['txt'] and Age (years) formula ---
$casesQuery = @"
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)

WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE  f.file_format IN ['txt']
return DISTINCT ss.study_subject_id as ``Case ID``,
   p.program_acronym as ``Program Code``,
    p.program_id as Program_ID,
   s.study_acronym as ``Arm``,
   ss.disease_subtype as ``Diagnosis``,
   sf.grouped_recurrence_score AS ``Recurrence Score``,
   d.tumor_size_group AS ``tumor_size``,
   d.er_status AS ``ER Status``,
   d.pr_status AS ``PR Status``,
   coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '') AS ``Age (years)``,
	demo.survival_time AS ``Survival (days)``
"@
$ws.Range("B2").Value = $casesQuery

# --- Count query text (reused verbatim in C3/C4) ---
$countQuery = @"
MATCH (ss:study_subject)
MATCH (samp:sample)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE f.file_format IN ["txt"]
WITH ss 
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
"@

# --- SamplesTab query text ---
$samplesQuery = @"
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE  f.file_format IN ['txt']
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS ``Sample ID``,
            ss.study_subject_id AS ``Case ID``,
            p.program_acronym AS ``Program Code``,
            s.study_acronym AS ``Arm``,
            ss.disease_subtype AS ``Diagnosis``,
            samp.tissue_type AS ``Tissue Type``,
            samp.composition AS ``Tissue Composition``,
            samp.sample_anatomic_site AS ``Sample Anatomic Site``,
            samp.method_of_sample_procurement AS ``Sample Procurement Method``
"@

# --- FilesTab query text ---
$filesQuery = @"
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE  f.file_format IN ['txt']
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS ``File Name``,
    head(labels(samp)) AS ``Association``,
    f.file_description AS ``Description``,
    f.file_format AS ``File Format``,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS ``Program Code``,
    s.study_acronym AS ``Arm``,
    ss.study_subject_id AS ``Case ID``,
    samp.sample_id AS ``Sample ID``
    order by f.file_name
"@

# Tab-name cells first (A3, A4) to mirror the shared-string allocation order
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# Row 3: SamplesTab
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $countQuery
$ws.Range("D3").Value = "TC01_Bento_Filter_FileType-text_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Bento_Filter_FileType-text_WebData.xlsx"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Row 4: FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $countQuery
$ws.Range("D4").Value = "TC01_Bento_Filter_FileType-text_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC01_Bento_Filter_FileType-text_WebData.xlsx"
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights to approximate Excel's auto-fit for the wrapped content
$ws.Range("B2").RowHeight = 403.2
$ws.Range("B3").RowHeight = 345.6
$ws.Range("B4").RowHeight = 409.6

# View state: selection + top-left cell as recorded after edit
$ws.Range("A4").Select()
$ws.Range("C2:E4").Select()
